$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Make room: the "Docentes responsaveis" block needs two data rows (for the
#    two professors) and everything below it shifts down by two rows.
#    Inserting whole rows preserves existing content/styles/heights for all
#    the rows that just move down.
# ---------------------------------------------------------------------------
$ws.Rows("13:14").Insert()

# ---------------------------------------------------------------------------
# Helper text blocks (kept as variables to keep the Range.Value assignments
# short and readable).
# ---------------------------------------------------------------------------
$objetivosPt   = "Introduzir o uso e prática de métodos numéricos para a solução de problemas matemáticos aplicados à Física e Engenharia. O estudante estará capacitado a descrever matematicamente e resolver numericamente problemas com o auxílio de algoritmos computacionais."
$prof1         = "3480026 - João Paulo Pascon"
$prof2         = "1176388 - Luiz Tadeu Fernandes Eleno"
$programaResPt = "Representação computacional de números em ponto flutuante; Zeros de funções; Sistemas de equações lineares; Método dos Mínimos Quadrados; Interpolação; Integração numérica; equações diferenciais ordinárias."
$programaPt    = "• Números em ponto flutuante: representação e precisão. • Raízes de funções: método da bissecção; método da falsa posição; método de Newton-Raphson; • Solução de sistemas de equações lineares: pivotamento e escalonamento; método de Gauss. • Método dos mínimos quadrados: ajuste de funções lineares nos parâmetros ajustáveis; ajuste de funções linearizáveis; ajuste de funções não-lineares usando a biblioteca scipy.optimize • Interpolação: método de Lagrange; método de Newton; • Integração numérica: regra dos trapézios; regra de Simpson; métodos avançados implementados na biblioteca scipy.integrate.  • Solução de equações diferenciais ordinárias: método de Euler; método de Runge-Kutta; métodos mais avançados da biblioteca scipy.integrate."
$metodo        = "Aulas expositivas e em laboratório computacional, trabalhos e exercícios comentados."
$criterio      = "Média aritmética de trabalhos propostos ao longo do curso (30%) e duas avaliações individuais (70%)."
$normaRecup    = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$bibliografia  = "Cunha, M. C. C., Métodos Numéricos. Editora Unicamp, 1993.Sperandio, D., Mendes, J. T., Monken e Silva, L. H. Cálculo Numérico. Pearson, 2003 LANGTANGEN, Hans Petter. A Primer on scientific programming with Python, 2a ed. New York: Springer, 2011. LANGTANGEN, Hans Petter. Python scripting for computational science, 5a ed. New York: Springer, 2016. SCOPATZ, A.; HUFF, K. D. Effective computation in physics: field guide to research in Python. Sebastpol, CA: O’Reilly Media, 2015."

# ---------------------------------------------------------------------------
# Small helper: give a brand-new (blank) cell the same look (wrap text,
# vertical alignment, font colour/weight, ...) as a known-good neighbour in
# the same column, then write its value. Using copy/paste-special for the
# formats (rather than poking individual Font/WrapText properties) makes the
# new cell share the exact same style record as its neighbour instead of
# Excel minting a brand new one.
# ---------------------------------------------------------------------------
function Set-StyledValue($targetAddr, $templateAddr, $value) {
    $ws.Range($templateAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($targetAddr).Value = $value
}

# ---------------------------------------------------------------------------
# 2) Objetivos / Objectives (row 10/11): these previously (wrongly) echoed
#    the professor's name - put the real objective text in place.
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = $objetivosPt
$ws.Range("C10").Value = $objetivosPt

# ---------------------------------------------------------------------------
# 3) Docentes responsaveis (rows 13 & 14, newly inserted): one professor per
#    row.
# ---------------------------------------------------------------------------
Set-StyledValue "B13" "B15" $prof1
Set-StyledValue "C13" "C15" $prof1
Set-StyledValue "B14" "B15" $prof2
Set-StyledValue "C14" "C15" $prof2

# ---------------------------------------------------------------------------
# 4) Programa resumido (row 15)
# ---------------------------------------------------------------------------
$ws.Range("B15").Value = $programaResPt
$ws.Range("C15").Value = $programaResPt

# ---------------------------------------------------------------------------
# 5) Programa (row 17)
# ---------------------------------------------------------------------------
$ws.Range("B17").Value = $programaPt
$ws.Range("C17").Value = $programaPt

# ---------------------------------------------------------------------------
# 6) Metodo (row 20)
# ---------------------------------------------------------------------------
$ws.Range("B20").Value = $metodo
$ws.Range("C20").Value = $metodo

# ---------------------------------------------------------------------------
# 7) Criterio (row 21)
# ---------------------------------------------------------------------------
$ws.Range("B21").Value = $criterio
$ws.Range("C21").Value = $criterio

# ---------------------------------------------------------------------------
# 8) Norma de recuperacao (row 22)
# ---------------------------------------------------------------------------
$ws.Range("B22").Value = $normaRecup
$ws.Range("C22").Value = $normaRecup

# ---------------------------------------------------------------------------
# 9) Bibliografia (row 23)
# ---------------------------------------------------------------------------
$ws.Range("B23").Value = $bibliografia
$ws.Range("C23").Value = $bibliografia

# ---------------------------------------------------------------------------
# 10) Column layout bug-fix: column A's width override used to bleed into
#     column B (min=1 max=2); column B has its own explicit width entry
#     right after it, so split the ranges so col A only covers column 1.
# ---------------------------------------------------------------------------
$ws.Columns(1).ColumnWidth = $ws.Columns(1).ColumnWidth
